$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1, matching the style of the existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data values for columns I and J, rows 2-13
$data = @{
    2  = @(8, 9)
    3  = @(5, 6)
    4  = @(11, 11)
    5  = @(8, 9)
    6  = @(9, 9)
    7  = @(5, 5)
    8  = @(7, 7)
    9  = @(6, 7)
    10 = @(9, 9)
    11 = @(5, 5)
    12 = @(3, 4)
    13 = @(3, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
